# MAJ documents avancement et tâches
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")
$chart = $ws.ChartObjects(1).Chart

# --- Fix a shared-string typo (double space) -----------------------------
$ws.Range("C12").Value = "Intégrer la nouvelle implémentation du joueur"

# --- Fill in the H/I progress columns for the existing tasks -------------
$ws.Range("H4").Value = 41428
$ws.Range("I4").Value = 41429

$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 0

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 6

$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 4

$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5

$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5

$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0

$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 5

# --- Move the chart's anchor one row down before the new row shifts ------
# everything beneath it (the engine does not auto-shift drawing anchors on
# a row insert, so this has to happen relative to the pre-insert layout).
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + 15

# --- Insert the new task "#9 - Création du menu" as row 13 ---------------
$ws.Rows.Item(13).Insert()
$ws.Range("B12:J12").Copy()
$ws.Range("B13:J13").PasteSpecial(-4122)

$ws.Range("B13").Value = "#9"
$ws.Range("C13").Value = "Création du menu"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 5

# --- Repair the totals row (old row 13, now row 14) -----------------------
$ws.Range("E14").Formula = '=SUM(E5:E13)'
$ws.Range("F14").Formula = '=SUM(F5:F13)'
$ws.Range("G14").Formula = '=SUM(G5:G13)'
$ws.Range("H14").Formula = '=SUM(H5:H13)'
$ws.Range("I14").Formula = '=SUM(I5:I13)'
$ws.Range("J14").Formula = '=SUM(J5:J13)'

# --- Repair the RAF% / date helper table (old rows 16-20, now 17-21) -----
$ws.Range("I17").Formula = '=F14/$E$14'
$ws.Range("J17").Formula = '=F$4'

$ws.Range("I18").Formula = '=G14/$E$14'
$ws.Range("J18").Formula = '=G4'

$ws.Range("I19").Formula = '=H14/$E$14'
$ws.Range("J19").Value = 41428

$ws.Range("I20").Formula = '=I14/$E$14'
$ws.Range("J20").Value = 41429

# --- Update the chart series to point at the new data range --------------
$ser = $chart.SeriesCollection(1)
$ser.Formula = '=SERIES(Sprint2!$J$17:$J$21,,Sprint2!$I$17:$I$21,1)'

# --- Print area grows by one row to keep the new task in range -----------
$ws.PageSetup.PrintArea = '$B$2:$K$34'

# --- Selection / dimension bookkeeping to match the saved file -----------
$ws.Range("I21").Select()
